# "Excelbestand levensonderhoud uitgebreid met query" -
# add a new worksheet "oef2_1" with a small people/characteristics table.

$xlCenter = -4108

$wb = $excel.ActiveWorkbook

# New sheet goes after the last existing tab (x-y, studietijd-cijfer, oef2_1).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "oef2_1"

# Headers (iq header is filled in later, after the numeric columns).
$ws.Range("A1").Value = "voornaam"
$ws.Range("B1").Value = "geslacht"
$ws.Range("C1").Value = "haarkleur"
$ws.Range("D1").Value = "lengte"
$ws.Range("E1").Value = "gewicht"

# Data entered column by column.
$ws.Range("A2").Value = "Chris"
$ws.Range("A3").Value = "Mari"
$ws.Range("A4").Value = "Otto"
$ws.Range("A5").Value = "Peter"
$ws.Range("A6").Value = "Vicky"

$ws.Range("B2").Value = "m"
$ws.Range("B3").Value = "v"
$ws.Range("B4").Value = "m"
$ws.Range("B5").Value = "m"
$ws.Range("B6").Value = "v"

$ws.Range("C2").Value = "bruin"
$ws.Range("C3").Value = "blond"
$ws.Range("C4").Value = "blond"
$ws.Range("C5").Value = "zwart"
$ws.Range("C6").Value = "rood"

$ws.Range("D2").Value = "groot"
$ws.Range("D3").Value = "groot"
$ws.Range("D4").Value = "normaal"
$ws.Range("D5").Value = "normaal"
$ws.Range("D6").Value = "klein"

$ws.Range("E2").Value = 185
$ws.Range("E3").Value = 176
$ws.Range("E4").Value = 181
$ws.Range("E5").Value = 178
$ws.Range("E6").Value = 164

$ws.Range("F2").Value = 95
$ws.Range("F3").Value = 104
$ws.Range("F4").Value = 98
$ws.Range("F5").Value = 108
$ws.Range("F6").Value = 112

# iq column header added last.
$ws.Range("F1").Value = "iq"

# Header row: bold + centered.
foreach ($cell in $ws.Range("A1:F1")) {
    $cell.HorizontalAlignment = $xlCenter
    $cell.Font.Bold = $true
}

# Data rows: centered.
$ws.Range("A2:F6").HorizontalAlignment = $xlCenter

# Leave the cursor below the table, like the author did.
[void]$ws.Range("A8").Select()
